$wb = $excel.ActiveWorkbook

# --- Sheet "Rushing" ---
$rushing = $wb.Worksheets.Item("Rushing")

# Row 4 - D.Cook
$rushing.Cells.Item(4, 3).Value = 155   # C4 1DATT
$rushing.Cells.Item(4, 4).Value = 110   # D4 2DATT
$rushing.Cells.Item(4, 5).Value = 25    # E4 3DATT

# Row 5 - A.Mattison
$rushing.Cells.Item(5, 3).Value = 41    # C5 1DATT
$rushing.Cells.Item(5, 4).Value = 27    # D5 2DATT

# Row 6 - K.Nwangwu
$rushing.Cells.Item(6, 5).Value = 2     # E6 3DATT

# Row 8 - J.Jefferson
$rushing.Cells.Item(8, 4).Value = 1     # D8 2DATT

# --- Sheet "Receiving" ---
$receiving = $wb.Worksheets.Item("Receiving")

# Row 2 - D.Cook
$receiving.Cells.Item(2, 3).Value = 48  # C2 Short Target
$receiving.Cells.Item(2, 4).Value = 33  # D2 Short Comp

# Row 3 - A.Mattison
$receiving.Cells.Item(3, 3).Value = 33  # C3 Short Target
$receiving.Cells.Item(3, 4).Value = 26  # D3 Short Comp

# Row 5 - D.Westbrook
$receiving.Cells.Item(5, 3).Value = 15  # C5 Short Target
$receiving.Cells.Item(5, 4).Value = 14  # D5 Short Comp

# Row 6 - C.Herndon
$receiving.Cells.Item(6, 3).Value = 118 # C6 Short Target
$receiving.Cells.Item(6, 4).Value = 82  # D6 Short Comp
$receiving.Cells.Item(6, 5).Value = 54  # E6 Deep Target
$receiving.Cells.Item(6, 6).Value = 28  # F6 Deep Comp

# Row 8 - B.Ellefson
$receiving.Cells.Item(8, 3).Value = 61  # C8 Short Target
$receiving.Cells.Item(8, 5).Value = 20  # E8 Deep Target
$receiving.Cells.Item(8, 6).Value = 10  # F8 Deep Comp

# Row 10 - I.Smith
$receiving.Cells.Item(10, 3).Value = 81 # C10 Short Target
$receiving.Cells.Item(10, 4).Value = 57 # D10 Short Comp

# Row 12 - L.Stocker
$receiving.Cells.Item(12, 3).Value = 3  # C12 Short Target
$receiving.Cells.Item(12, 4).Value = 3  # D12 Short Comp
$receiving.Cells.Item(12, 5).Value = 3  # E12 Deep Target
$receiving.Cells.Item(12, 6).Value = 2  # F12 Deep Comp
